$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.3227736666666667
$ws.Range("H2").Value = 0.968321
$ws.Range("I2").Value = 0.1416094457286952
$ws.Range("J2").Value = 0.1416094457286952
$ws.Range("M2").Value = 4.066173333333333
$ws.Range("N2").Value = 12.19852
$ws.Range("O2").Value = 0.8070107842953054
$ws.Range("P2").Value = 0.8070107842953055
$ws.Range("Q2").Value = 1.312453676102222
$ws.Range("R2").Value = 11.81208308492
$ws.Range("S2").Value = 0.1142803498611378
$ws.Range("T2").Value = 0.1142803498611378

# Row 3
$ws.Range("G3").Value = 0.3227736666666667
$ws.Range("H3").Value = 0.968321
$ws.Range("I3").Value = 0.1416094457286952
$ws.Range("J3").Value = 0.1416094457286952
$ws.Range("O3").Value = 0.03600419273120554
$ws.Range("P3").Value = 0.03600419273120554
$ws.Range("Q3").Value = 0.05855415568755555
$ws.Range("R3").Value = 0.5269874011879999
$ws.Range("S3").Value = 0.005098533776575132
$ws.Range("T3").Value = 0.005098533776575132

# Row 4
$ws.Range("G4").Value = 0.3227736666666667
$ws.Range("H4").Value = 0.968321
$ws.Range("I4").Value = 0.1416094457286952
$ws.Range("J4").Value = 0.1416094457286952
$ws.Range("M4").Value = 0.7909786666666667
$ws.Range("N4").Value = 2.372936
$ws.Range("O4").Value = 0.156985022973489
$ws.Range("P4").Value = 0.156985022973489
$ws.Range("Q4").Value = 0.2553070844951111
$ws.Range("R4").Value = 2.297763760456
$ws.Range("S4").Value = 0.02223056209098226
$ws.Range("T4").Value = 0.02223056209098226

# Row 5
$ws.Range("I5").Value = 0.8226066833587575
$ws.Range("J5").Value = 0.8226066833587576
$ws.Range("M5").Value = 4.066173333333333
$ws.Range("N5").Value = 12.19852
$ws.Range("O5").Value = 0.8070107842953054
$ws.Range("P5").Value = 0.8070107842953055
$ws.Range("Q5").Value = 7.624019428964443
$ws.Range("R5").Value = 68.61617486067999
$ws.Range("S5").Value = 0.6638524647039108
$ws.Range("T5").Value = 0.663852464703911

# Row 6
$ws.Range("I6").Value = 0.8226066833587575
$ws.Range("J6").Value = 0.8226066833587576
$ws.Range("O6").Value = 0.03600419273120554
$ws.Range("P6").Value = 0.03600419273120554
$ws.Range("S6").Value = 0.02961728956962648
$ws.Range("T6").Value = 0.02961728956962648

# Row 7
$ws.Range("I7").Value = 0.8226066833587575
$ws.Range("J7").Value = 0.8226066833587576
$ws.Range("M7").Value = 0.7909786666666667
$ws.Range("N7").Value = 2.372936
$ws.Range("O7").Value = 0.156985022973489
$ws.Range("P7").Value = 0.156985022973489
$ws.Range("Q7").Value = 1.483074189958222
$ws.Range("R7").Value = 13.347667709624
$ws.Range("S7").Value = 0.1291369290852201
$ws.Range("T7").Value = 0.1291369290852202

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.081563
$ws.Range("H8").Value = 0.244689
$ws.Range("I8").Value = 0.03578387091254728
$ws.Range("J8").Value = 0.03578387091254728
$ws.Range("M8").Value = 4.066173333333333
$ws.Range("N8").Value = 12.19852
$ws.Range("O8").Value = 0.8070107842953054
$ws.Range("P8").Value = 0.8070107842953055
$ws.Range("Q8").Value = 0.3316492955866666
$ws.Range("R8").Value = 2.98484366028
$ws.Range("S8").Value = 0.02887796973025674
$ws.Range("T8").Value = 0.02887796973025675

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.081563
$ws.Range("H9").Value = 0.244689
$ws.Range("I9").Value = 0.03578387091254728
$ws.Range("J9").Value = 0.03578387091254728
$ws.Range("O9").Value = 0.03600419273120554
$ws.Range("P9").Value = 0.03600419273120554
$ws.Range("Q9").Value = 0.01479628945466666
$ws.Range("R9").Value = 0.133166605092
$ws.Range("S9").Value = 0.001288369385003932
$ws.Range("T9").Value = 0.001288369385003932

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.081563
$ws.Range("H10").Value = 0.244689
$ws.Range("I10").Value = 0.03578387091254728
$ws.Range("J10").Value = 0.03578387091254728
$ws.Range("M10").Value = 0.7909786666666667
$ws.Range("N10").Value = 2.372936
$ws.Range("O10").Value = 0.156985022973489
$ws.Range("P10").Value = 0.156985022973489
$ws.Range("Q10").Value = 0.06451459298933333
$ws.Range("R10").Value = 0.580631336904
$ws.Range("S10").Value = 0.0056175317972866
$ws.Range("T10").Value = 0.0056175317972866
